# Update column C ("Förändrad") date value for rows 2-13
# from 2023-10-22 (serial 45221) to 2023-10-25 (serial 45224)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 13; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45221) {
        $cell.Value2 = 45224
    }
}
